# "make PRS analysis slide white"
# Slide 11 ("PRS analysis"): give it an explicit white (bg1) background,
# and reposition/resize the picture + the p<=0.0024% textbox to match.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

# --- Add an explicit slide background: solid fill, scheme color bg1 (white) ---
$s.Background.Fill.Solid()
$s.Background.Fill.ForeColor.SchemeColor = "bg1"

# --- Reposition / resize the picture (Picture 3) ---
$pic = $s.Shapes.Item(2)
$emuPerPt = 914400 / 72
$pic.Left   = 1012825 / $emuPerPt
$pic.Top    = -259557 / $emuPerPt
$pic.Width  = 9836150 / $emuPerPt
$pic.Height = 7377113 / $emuPerPt

# --- Reposition the "p<=0.0024%" textbox (TextBox 6); size unchanged ---
$txt = $s.Shapes.Item(3)
$txt.Left = 2425700 / $emuPerPt
$txt.Top  = 5373296 / $emuPerPt
